$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new row (professor name) was inserted at row 13, pushing the rest of the
# table (old rows 13-23) down to rows 14-24.
$ws.Rows("13:13").Insert()

# The insert carries column-A formatting into the new row; row 13 has no
# label cell in the target layout, so drop it.
$ws.Range("A13").Clear()

# Row 10 (Objetivos): replace placeholder text with the real objectives text
$ws.Range("B10").Value = 'Apresentar aos alunos o panorama geral da administração estratégica de marketing, capacitando-os a atuar no processo gerencial de marketing sob as perspectivas estratégica e operacional.'
$ws.Range("C10").Value = 'Apresentar aos alunos o panorama geral da administração estratégica de marketing, capacitando-os a atuar no processo gerencial de marketing sob as perspectivas estratégica e operacional.'

# Row 13 (Docentes responsaveis value): professor name
$ws.Range("B13").Value = '11079086 - Herlandí de Souza Andrade'
$ws.Range("C13").Value = '11079086 - Herlandí de Souza Andrade'
# Match the wrap-top (plain/red) formatting used by the rest of the table
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 14 (Programa resumido): new short-syllabus text (PT)
$ws.Range("B14").Value = 'Marketing estratégico e planejamento estratégico orientado para o mercado.'
$ws.Range("C14").Value = 'Marketing estratégico e planejamento estratégico orientado para o mercado.'

# Row 16 (Programa): new detailed PT syllabus text
$ws.Range("B16").Value = 'MARKETING ESTRATÉGICO1. Evolução do Conceito de Marketing e Sistema de Marketing2. Marketing, Conceito de Valor, Orientação para Mercado3. Análise de Mercado e Comportamento do Consumidor4. Sistema de Informações de Marketing e Inteligência de Mercado5. Modalidades de Marketing6. Fundamentos de Estratégia Empresarial e Marketing Estratégico7. Administração Estratégica e Marketing Estratégico8. Instrumentos Analíticos para Avaliar Oportunidades de Mercado9. Segmentação de Mercado e Posicionamento10. O Plano Estratégico de Marketing11. Comunicação: assessorias de imprensa, SAC''S, Ombudsman'
$ws.Range("C16").Value = 'MARKETING ESTRATÉGICO1. Evolução do Conceito de Marketing e Sistema de Marketing2. Marketing, Conceito de Valor, Orientação para Mercado3. Análise de Mercado e Comportamento do Consumidor4. Sistema de Informações de Marketing e Inteligência de Mercado5. Modalidades de Marketing6. Fundamentos de Estratégia Empresarial e Marketing Estratégico7. Administração Estratégica e Marketing Estratégico8. Instrumentos Analíticos para Avaliar Oportunidades de Mercado9. Segmentação de Mercado e Posicionamento10. O Plano Estratégico de Marketing11. Comunicação: assessorias de imprensa, SAC''S, Ombudsman'

# Row 19 (Metodo): teaching method text
$ws.Range("B19").Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'
$ws.Range("C19").Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'

# Row 20 (Criterio): grading criteria text
$ws.Range("B20").Value = 'Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas.'
$ws.Range("C20").Value = 'Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas.'

# Row 21 (Norma de recuperacao): makeup exam rule text
$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'

# Row 22 (Bibliografia): bibliography text
$ws.Range("B22").Value = 'KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.KOTLER, P.; KARTAJAYA, H.; SETIAWAN, I. Marketing 4.0: do Tradicional ao Digital. São Paulo: Sextante, 2017.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L.  Marketing Essencial. 5 ed. São Paulo: Pearson, 2013.SANDHUSEN, R. L. Marketing Básico - Série Essencial. 3 ed. São Paulo: Saraiva, 2010.SAPIRO, Arão., CHIAVENATO, I. Planejamento Estratégico. Campus, 2ª. edição, 2010 KOTLER, P. Administração de Marketing, edição do milênio, revisão técnica de Prof. Arão Sapiro. Prentice-Hall, 2000. HOOLEY, Graham J.; PIERCY, Nigel F.; SAUNDERS, John A. Estratégia de Marketing e Posicionamento Competitivo tradução e revisão técnica: Prof. Arão Sapiro. Pearson Education do Brasil, 2001. SAPIRO, ARAO; GANGANA, MAURÍCIO; LIMA, MIGUEL; VILHENA, JOÃO BAPTISTA. Gestão de Marketing . FGV Editora, 2004. BOONE, L. e KURTZ, D.L. Marketing contemporâneo. 8ª ed. São Paulo, Livros Técnicos e Científicos, 1998. KOTLER, P; JATURISPITAK, S. e MAESINCIE, S. O marketing das nações. São Paulo, Futura, 1997. MARTINS, J.R. e BLECHER, N. O império das marcas. 2ª ed. São Paulo, Negócio Editora, 1997 THUROW, L.C. O futuro do capitalismo. 2ª ed. São Paulo, Rocco, 1997. VAZ, G. N. Marketing institucional. São Paulo, Pioneira, 1995. Bibliografia Complementar Artigos das Revistas: Marketing, Meio e Mensagem, Exame, Dinheiro, Revista da Escola de Administração da FEA-USP, Revista ESPM.'
$ws.Range("C22").Value = 'KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.KOTLER, P.; KARTAJAYA, H.; SETIAWAN, I. Marketing 4.0: do Tradicional ao Digital. São Paulo: Sextante, 2017.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L.  Marketing Essencial. 5 ed. São Paulo: Pearson, 2013.SANDHUSEN, R. L. Marketing Básico - Série Essencial. 3 ed. São Paulo: Saraiva, 2010.SAPIRO, Arão., CHIAVENATO, I. Planejamento Estratégico. Campus, 2ª. edição, 2010 KOTLER, P. Administração de Marketing, edição do milênio, revisão técnica de Prof. Arão Sapiro. Prentice-Hall, 2000. HOOLEY, Graham J.; PIERCY, Nigel F.; SAUNDERS, John A. Estratégia de Marketing e Posicionamento Competitivo tradução e revisão técnica: Prof. Arão Sapiro. Pearson Education do Brasil, 2001. SAPIRO, ARAO; GANGANA, MAURÍCIO; LIMA, MIGUEL; VILHENA, JOÃO BAPTISTA. Gestão de Marketing . FGV Editora, 2004. BOONE, L. e KURTZ, D.L. Marketing contemporâneo. 8ª ed. São Paulo, Livros Técnicos e Científicos, 1998. KOTLER, P; JATURISPITAK, S. e MAESINCIE, S. O marketing das nações. São Paulo, Futura, 1997. MARTINS, J.R. e BLECHER, N. O império das marcas. 2ª ed. São Paulo, Negócio Editora, 1997 THUROW, L.C. O futuro do capitalismo. 2ª ed. São Paulo, Rocco, 1997. VAZ, G. N. Marketing institucional. São Paulo, Pioneira, 1995. Bibliografia Complementar Artigos das Revistas: Marketing, Meio e Mensagem, Exame, Dinheiro, Revista da Escola de Administração da FEA-USP, Revista ESPM.'

